{"js": "// Update the multiplication-table answer cells with newly generated values.\n// Each original \"A\u00d7B=C\" string is unique in the document, so searching the\n// body for each literal string and replacing the single match is sufficient\n// to retarget exactly the cells the diff touches.\n\nconst replacements = [\n  { find: \"20\u00d722=440\", replace: \"93\u00d750=4650\" },\n  { find: \"67\u00d791=6097\", replace: \"12\u00d729=348\" },\n  { find: \"20\u00d796=1920\", replace: \"58\u00d790=5220\" },\n  { find: \"29\u00d725=725\", replace: \"13\u00d741=533\" },\n  { find: \"98\u00d714=1372\", replace: \"54\u00d713=702\" },\n  { find: \"20\u00d785=1700\", replace: \"90\u00d724=2160\" },\n  { find: \"67\u00d724=1608\", replace: \"74\u00d718=1332\" },\n  { find: \"15\u00d795=1425\", replace: \"52\u00d759=3068\" },\n  { find: \"66\u00d711=726\", replace: \"94\u00d740=3760\" },\n  { find: \"47\u00d716=752\", replace: \"15\u00d741=615\" },\n  { find: \"28\u00d730=840\", replace: \"33\u00d784=2772\" },\n  { find: \"21\u00d712=252\", replace: \"82\u00d718=1476\" },\n  { find: \"69\u00d739=2691\", replace: \"40\u00d760=2400\" },\n  { find: \"79\u00d757=4503\", replace: \"36\u00d754=1944\" },\n  { find: \"65\u00d768=4420\", replace: \"89\u00d784=7476\" },\n  { find: \"38\u00d763=2394\", replace: \"19\u00d753=1007\" },\n  { find: \"36\u00d713=468\", replace: \"11\u00d718=198\" },\n  { find: \"33\u00d759=1947\", replace: \"63\u00d766=4158\" },\n  { find: \"53\u00d767=3551\", replace: \"82\u00d747=3854\" },\n  { find: \"61\u00d798=5978\", replace: \"89\u00d751=4539\" },\n  { find: \"91\u00d736=3276\", replace: \"74\u00d718=1332\" },\n  { find: \"56\u00d717=952\", replace: \"85\u00d746=3910\" },\n  { find: \"25\u00d743=1075\", replace: \"28\u00d783=2324\" },\n  { find: \"32\u00d751=1632\", replace: \"17\u00d767=1139\" },\n  { find: \"14\u00d738=532\", replace: \"18\u00d778=1404\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-table answer cells with newly generated values.\n# Each original \"A\u00d7B=C\" string is unique in the document, so a sequence of\n# literal (non-wildcard) Find/Replace operations over the whole document\n# content safely retargets exactly the cells the diff touches.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Find=\"20\u00d722=440\";  Replace=\"93\u00d750=4650\"},\n    @{Find=\"67\u00d791=6097\"; Replace=\"12\u00d729=348\"},\n    @{Find=\"20\u00d796=1920\"; Replace=\"58\u00d790=5220\"},\n    @{Find=\"29\u00d725=725\";  Replace=\"13\u00d741=533\"},\n    @{Find=\"98\u00d714=1372\"; Replace=\"54\u00d713=702\"},\n    @{Find=\"20\u00d785=1700\"; Replace=\"90\u00d724=2160\"},\n    @{Find=\"67\u00d724=1608\"; Replace=\"74\u00d718=1332\"},\n    @{Find=\"15\u00d795=1425\"; Replace=\"52\u00d759=3068\"},\n    @{Find=\"66\u00d711=726\";  Replace=\"94\u00d740=3760\"},\n    @{Find=\"47\u00d716=752\";  Replace=\"15\u00d741=615\"},\n    @{Find=\"28\u00d730=840\";  Replace=\"33\u00d784=2772\"},\n    @{Find=\"21\u00d712=252\";  Replace=\"82\u00d718=1476\"},\n    @{Find=\"69\u00d739=2691\"; Replace=\"40\u00d760=2400\"},\n    @{Find=\"79\u00d757=4503\"; Replace=\"36\u00d754=1944\"},\n    @{Find=\"65\u00d768=4420\"; Replace=\"89\u00d784=7476\"},\n    @{Find=\"38\u00d763=2394\"; Replace=\"19\u00d753=1007\"},\n    @{Find=\"36\u00d713=468\";  Replace=\"11\u00d718=198\"},\n    @{Find=\"33\u00d759=1947\"; Replace=\"63\u00d766=4158\"},\n    @{Find=\"53\u00d767=3551\"; Replace=\"82\u00d747=3854\"},\n    @{Find=\"61\u00d798=5978\"; Replace=\"89\u00d751=4539\"},\n    @{Find=\"91\u00d736=3276\"; Replace=\"74\u00d718=1332\"},\n    @{Find=\"56\u00d717=952\";  Replace=\"85\u00d746=3910\"},\n    @{Find=\"25\u00d743=1075\"; Replace=\"28\u00d783=2324\"},\n    @{Find=\"32\u00d751=1632\"; Replace=\"17\u00d767=1139\"},\n    @{Find=\"14\u00d738=532\";  Replace=\"18\u00d778=1404\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
